# Applies the commit's changes to the "Avverkningsanmälningar" sheet:
#   1. Column C ("Förändrad") date serial 45184 -> 45186 for every data row.
#   2. The HYPERLINK() formulas in columns S, T, V, W, X, Y (rows 2-17, the
#      rows that actually have species findings / links) gain a second
#      HYPERLINK argument: the friendly display text, which is the row's
#      "Beteckning" (column A) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $designation = $ws.Cells.Item($r, 1).Value2

    if ($designation -eq $null -or $designation -eq "") {
        continue
    }

    # 1) Bump the "Förändrad" date serial from 45184 to 45186.
    $cDate = $ws.Cells.Item($r, 3).Value2
    if ($cDate -eq 45184) {
        $ws.Cells.Item($r, 3).Value2 = 45186
    }

    # 2) Add the friendly-name second argument to each HYPERLINK formula
    #    present on this row (columns S, T, V, W, X, Y = 19, 20, 22, 23, 24, 25).
    foreach ($col in 19, 20, 22, 23, 24, 25) {
        $cell = $ws.Cells.Item($r, $col)
        $formula = $cell.Formula
        if ([string]::IsNullOrEmpty($formula)) {
            continue
        }
        if ($formula -like '*HYPERLINK(*' -and $formula -notlike '*,*') {
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $designation + '")'
            $cell.Formula = $newFormula
        }
    }
}
